# issue #5: property land done
# Clean up stray internal spaces / OCR typos in shared-string text across
# the "存款" (deposit), "基金受益憑證" (fund), "保險" (insurance) and
# "債務" (debt) sheets, and drop the duplicate "新臺幣'" string in favour
# of the existing "新臺幣".

$wb = $excel.ActiveWorkbook

# -------- 存款 (deposits) --------
$ws = $wb.Worksheets.Item("存款")
$ws.Range("B2").Value  = "台北富邦商業銀行城中分行"
$ws.Range("B3").Value  = "台北富邦商業銀行城中分行"
$ws.Range("D8").Value  = "新臺幣"
$ws.Range("B9").Value  = "中華郵政股份有限公司高雄新庄仔郵局"
# "83，385" -> "83385": keep this a text cell (it's a shared string in the
# workbook, not a number) even though it now looks numeric.
$ws.Range("G9").Value  = "'83385"
$ws.Range("G9").Style  = "Normal"
$ws.Range("B21").Value = "花旗（台灣)商業銀行臺北分行"
$ws.Range("B22").Value = "中華郵政股份有限公司高雄新庄仔郵局"
$ws.Range("B23").Value = "中華郵政股份有限公司高雄新庄仔郵局"
$ws.Range("B24").Value = "花旗（台灣)商業銀行臺北分行"
$ws.Range("B26").Value = "兆豐國際商業銀行新竹分行"
$ws.Range("B27").Value = "中華郵政股份有限公司郵政儲金匯業局"

# -------- 基金受益憑證 (funds) --------
$ws = $wb.Worksheets.Item("基金受益憑證")
$ws.Range("B2").Value = "台灣工銀大眾基金"
$ws.Range("D2").Value = "台灣工銀證券投信公司"
$ws.Range("D3").Value = "國泰證券投信公司"

# -------- 保險 (insurance) --------
$ws = $wb.Worksheets.Item("保險")
$ws.Range("C2").Value = "南山金福利21年期還本養老壽險"
$ws.Range("C5").Value = "全球人壽全球104終身壽險甲型"
$ws.Range("C6").Value = "全球人壽全球104終身壽險甲型"
$ws.Range("C7").Value = "全球人壽全球104終身壽險甲塑"
$ws.Range("C8").Value = "全球人壽全球104終身壽險甲型"

# -------- 債務 (debt) --------
$ws = $wb.Worksheets.Item("債務")
$ws.Range("D2").Value = "大眾商業銀行新生分行臺北市中正區忠孝東路"
$ws.Range("F2").Value = "94年01月06日"
$ws.Range("D3").Value = "幸福人壽保險股份有限公司臺北市中正區忠孝西路"
$ws.Range("F3").Value = "95年12月25日"
$ws.Range("G3").Value = "般貸款"
$ws.Range("D4").Value = "幸福人壽保險股份有限公司臺北市中正區忠孝西路"
$ws.Range("F4").Value = "96年04月12日"
$ws.Range("G4").Value = "般貸款"
$ws.Range("D5").Value = "幸福人壽保險股份有限公司臺北市中正區忠孝西路"
$ws.Range("F5").Value = "96年05月28日"
$ws.Range("G5").Value = "般貸款"
$ws.Range("D6").Value = "幸福人壽保險股份有限公司臺北市中正區忠孝西路"
$ws.Range("F6").Value = "96年12月11曰"
